$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.293.42'
$ws.Range("E2").Value = '  -3.02%  '
$ws.Range("D3").Value = '''1.729.57'
$ws.Range("E3").Value = '  -3.86%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''321.72'
$ws.Range("E5").Value = '  -4.74%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("D7").Value = '''0.4227'
$ws.Range("E7").Value = '  -10.51%  '
$ws.Range("D8").Value = '''0.3577'
$ws.Range("E8").Value = '  -3.85%  '
$ws.Range("D9").Value = '''44.82'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").Value = '''0.07406'
$ws.Range("E10").Value = '  -3.63%  '
$ws.Range("D11").Value = '''1.106'
$ws.Range("D12").Value = '''1.003'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '''21.33'
$ws.Range("E13").Value = '  -5.93%  '
$ws.Range("D14").Value = '''6.050'
$ws.Range("E14").Value = '  -5.16%  '
$ws.Range("D15").Value = '''7.090'
$ws.Range("E15").Value = '  -4.27%  '
$ws.Range("D16").Value = '''1.729.35'
$ws.Range("E16").Value = '  -3.83%  '
$ws.Range("D17").Value = '''0.00001057'
$ws.Range("E17").Value = '  -3.65%  '
$ws.Range("D18").Value = '''86.54'
$ws.Range("E18").Value = '  +4.81%  '
$ws.Range("D19").Value = '''0.05948'
$ws.Range("E19").Value = '  -11.95%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Value = '''16.68'
$ws.Range("E21").Value = '  -4.42%  '
$ws.Range("D22").Value = '''6.058'
$ws.Range("E22").Value = '  -5.70%  '
$ws.Range("D23").Value = '''0.5234'
$ws.Range("E23").Value = '  -5.32%  '
$ws.Range("D24").Value = '''27.335.35'
$ws.Range("E24").Value = '  -2.89%  '
$ws.Range("D25").Value = '''11.29'
$ws.Range("E25").Value = '  -5.24%  '
$ws.Range("D26").Value = '''2.397'
$ws.Range("E26").Value = '  -0.46%  '
$ws.Range("D27").Value = '''20.06'
$ws.Range("E27").Value = '  -3.89%  '
$ws.Range("D28").Value = '''2.340'
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").Value = '''148.74'
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("D30").Value = '''1.925.82'
$ws.Range("E30").Value = '  -3.90%  '
$ws.Range("D31").Value = '''1.222'
$ws.Range("E31").Value = '  -3.25%  '
$ws.Range("D32").Value = '''125.76'
$ws.Range("E32").Value = '  -6.06%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = '''0.09066'
$ws.Range("E33").Value = '  -6.49%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''5.584'
$ws.Range("E34").Value = '  -6.03%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''3.613'
$ws.Range("E35").Value = '  -10.65%  '
$ws.Range("D36").Value = '''12.58'
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("D37").Value = '''0.2154'
$ws.Range("E37").Value = '  -2.95%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '''0.06093'
$ws.Range("E38").Value = '  -4.23%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.02247'
$ws.Range("E39").Value = '  -5.72%  '
$ws.Range("D40").Value = '''5.025'
$ws.Range("E40").Value = '  -4.70%  '
$ws.Range("D41").Value = '''0.6339'
$ws.Range("E41").Value = '  -5.61%  '
$ws.Range("D42").Value = '''1.182'
$ws.Range("E42").Value = '  -4.61%  '
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = '''1.411'
$ws.Range("E44").Value = '  -6.04%  '
$ws.Range("D45").Value = '''7.868'
$ws.Range("E45").Value = '  -2.77%  '
$ws.Range("D46").Value = '''13.48'
$ws.Range("E46").Value = '  -4.58%  '
$ws.Range("D47").Value = '''3.724'
$ws.Range("E47").Value = '  -3.38%  '
$ws.Range("D48").Value = '''0.5797'
$ws.Range("E48").Value = '  -6.14%  '
$ws.Range("D49").Value = '''124.51'
$ws.Range("E49").Value = '  -4.60%  '
$ws.Range("D50").Value = '''1.933'
$ws.Range("E50").Value = '  -6.38%  '
$ws.Range("D51").Value = '''0.06809'
$ws.Range("E51").Value = '  -4.48%  '